# Add "Bài 20" section (row 22) with its sorting-feature hyperlink entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$commitUrl = "https://github.com/nguyentienminh07102004/product-management/commit/d0a373564e97957e32fc6716531b8809e75cba35"
$label = "1. Sắp xếp sản phẩm theo các tiêu chí khác nhau"

$ws.Range("A22").Value = "Bài 20"
$ws.Range("B22").Value = $label

$ws.Hyperlinks.Add(
    $ws.Range("B22"),
    $commitUrl,
    [Type]::Missing,
    [Type]::Missing,
    $commitUrl
) | Out-Null

# Match the existing "Hyperlink" cell style used by the other link cells,
# and restore the display text Hyperlinks.Add may have overwritten.
$ws.Range("B22").Style = "Hyperlink"
$ws.Range("B22").Value = $label

$ws.Range("B22").Select()
